# "Fix some details for packer image builds"
#
# In the Cloud-Init template creation section, the parenthetical note
#   " (this file is in directory Ubuntu-2004-cloudinit-templates)"
# needs to become
#   " (this file is in directory Ubuntu-2004-cloudinit-template)"
# (the directory name is singular "template", not "templates"), and the
# run is split into several pieces around the (relocated) "_GoBack"
# bookmark, matching how Word naturally re-chunks text runs around an
# edit point.

$d = $word.ActiveDocument

# Locate the full parenthetical phrase so we do not depend on hard-coded
# character offsets.
$rngFull = $d.Content
$rngFull.Find.ClearFormatting()
$foundFull = $rngFull.Find.Execute(" (this file is in directory Ubuntu-2004-cloudinit-templates)", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundFull) {
    throw "Could not find the Ubuntu-2004-cloudinit-templates phrase"
}
$phraseStart = $rngFull.Start

# Work out (relative to the match) where the run boundaries need to land:
#  - right before " Ubuntu-2004-cloudinit-template..."
#  - right before the trailing "s)" (the "s" we are deleting)
$phrase = " (this file is in directory Ubuntu-2004-cloudinit-templates)"
$splitDirUbuntu = $phraseStart + $phrase.IndexOf(" Ubuntu")
$splitTemplateS = $phraseStart + $phrase.IndexOf("s)")

# Drop temporary bookmarks at each desired run boundary so the upcoming
# text edit cannot coalesce runs across them.
$d.Bookmarks.Add("ZZZ_barrier1", $d.Range($phraseStart, $phraseStart))
$d.Bookmarks.Add("ZZZ_barrier2", $d.Range($splitDirUbuntu, $splitDirUbuntu))
$d.Bookmarks.Add("ZZZ_barrier3", $d.Range($splitTemplateS, $splitTemplateS))

# Delete the "s" that turns "templates" into "template".
$rngS = $d.Range($splitTemplateS, $splitTemplateS + 1)
$rngS.Text = ""

# The document's "_GoBack" bookmark (last-edit marker) now belongs right
# after "template", i.e. exactly where we just edited - re-seat it there.
$d.Bookmarks.Add("_GoBack", $d.Range($splitTemplateS, $splitTemplateS))

# Clean up the helper bookmarks now that the run split has taken effect.
$d.Bookmarks.Item("ZZZ_barrier1").Delete()
$d.Bookmarks.Item("ZZZ_barrier2").Delete()
$d.Bookmarks.Item("ZZZ_barrier3").Delete()
